# Refatorando o consolidador para modelo ETL
# Overwrite the absenteeism data rows (2-11) with the new ETL-produced values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(36445, "Dr. João Lucas Costa", "Financeiro", "Problemas pessoais", 1, 45080, 11561.51),
    @(48955, "Lorena Barbosa",       "Engenharia", "Viagem de negócios", 5, 45099, 5320.15),
    @(72679, "Eloah Aragão",         "P&D",         "Doença",             3, 45079, 3564.26),
    @(54039, "Theo da Rocha",        "TI",          "Outros",             5, 45096, 3536.23),
    @(67407, "Brenda Porto",         "TI",          "Outros",             8, 45078, 5256.35),
    @(14002, "Raquel Costela",       "Jurídico",    "Problemas pessoais", 8, 45099, 12301.18),
    @(71002, "Bruna Rodrigues",      "TI",          "Consulta médica",    3, 45093, 12236.72),
    @(21226, "Francisco da Cruz",    "Operações",   "Outros",             7, 45102, 10316.83),
    @(33907, "Milena da Cruz",       "Operações",   "Viagem de negócios", 8, 45080, 8345.17),
    @(18544, "João Vitor Costa",     "Engenharia",  "Outros",             6, 45079, 6655.68)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
}
